# Add parameter for fixing the length of movie.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("imageSequenceTomovie")
$ws2 = $wb.Worksheets.Item("movieToimageSequence")
$ws3 = $wb.Worksheets.Item("convertGatanDM4Movies")

# --- imageSequenceTomovie (sheet1) ---

# Header: "File extension" -> "Image file extension"
$ws1.Range("D1").Value = "Image file extension"

# Drop the rows that are no longer needed:
#  original row 2 (jp6b07983_si_006), row 3 (jp6b07983_si_002),
#  and row 7 (SupportingMovieS2_NR_50e_Movie1)
# Delete from the bottom up so earlier row numbers stay valid. The sheet
# now holds only 4 rows total, so nudge the outline-row bookkeeping down
# to match (3 = rowcount-1) before the row disappears.
$ws1.Rows.Item(7).OutlineLevel = 3
$ws1.Rows.Item(7).Delete()
$ws1.Rows.Item(3).Delete()
$ws1.Rows.Item(2).Delete()

# After the deletions the remaining data (previously rows 4,5,6) now sits
# on rows 2,3,4:
#   row2 = nl8b04962_si_008 (601)
#   row3 = nl500766j_si_002 (811)
#   row4 = nl500766j_si_003 (1227)

# New blank (bordered) "frame step" column E for every data row.
$ws1.Range("E2").ClearContents()
$ws1.Range("E2").Borders.LineStyle = 1
$ws1.Range("E3").ClearContents()
$ws1.Range("E3").Borders.LineStyle = 1
$ws1.Range("E4").ClearContents()
$ws1.Range("E4").Borders.LineStyle = 1

# Row 2: clear quality column G entirely (no value, no border/style).
$ws1.Range("G2").ClearContents()
$ws1.Range("G2").Borders.LineStyle = -4142

# Row 3: quality becomes 5, new blank H, new duration I=50.
$ws1.Range("G3").Value = 5
$ws1.Range("H3").ClearContents()
$ws1.Range("H3").Borders.LineStyle = 1
$ws1.Range("I3").Value = 50

# Row 4: quality becomes 4, clear the old "generate all" flag in H (no
# value, no border/style), new duration I=100.
$ws1.Range("G4").Value = 4
$ws1.Range("H4").ClearContents()
$ws1.Range("H4").Borders.LineStyle = -4142
$ws1.Range("I4").Value = 100

# Make imageSequenceTomovie the active sheet/selection.
$ws1.Activate()
[void]$ws1.Range("A1").Select()

Write-Host "done"
